# "10Th - MB for single stock and added new group"
#
# The weekly MarketBeat analyst-rank sheet gets a new "10th" snapshot:
#   - two new date columns (Jun_27, Jun_26 x2) are inserted right after the
#     firm-name column, pushing the existing Jun_17/Jun_15/Jun_13/Jun_10
#     columns three slots to the right
#   - the new columns are seeded with the default "UN" rating for every
#     existing analyst row
#   - two new analyst rows (Benchmark, Evercore ISI) are appended at the
#     bottom of the table - a new coverage "group" that only has data for
#     the first few (most recent) date columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns right after column A (old B:E -> new E:H)
$ws.Range("B1:D1").EntireColumn.Insert()

# Header row for the 3 freshly inserted date columns (C/D share "Jun_26")
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Default "UN" rating for every existing analyst row in the new columns
$ws.Range("B2:D27").Value = "UN"

# New analyst/broker group appended at the bottom of the table
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# Most recent column filled in last
$ws.Range("B1").Value = "Jun_27"
